# Reverse the order of comma-separated "Recorded By" entries in column G.
# Cells with a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        $reversed = $parts[($parts.Count - 1)..0]
        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
